$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Address, $Text) {
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "52.290.82"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.011.16"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws "D5" "355.80"
Set-TextValue $ws "D6" "108.48"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue $ws "D9" "0.621"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("E11").Value = "  +2.00%  "
Set-TextValue $ws "D12" "0.0863"
$ws.Range("E12").Value = "  -3.95%  "
Set-TextValue $ws "D13" "19.34"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").Value = "3.489.29"
$ws.Range("E14").Value = "  +2.08%  "
Set-TextValue $ws "D15" "7.71"
$ws.Range("E15").Value = "  -4.35%  "
$ws.Range("D16").Value = "3.015.49"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "52.342.88"
Set-TextValue $ws "D19" "3.55"
$ws.Range("E19").Value = "  +8.52%  "
$ws.Range("E20").Value = "  -2.00%  "
Set-TextValue $ws "D21" "13.74"
$ws.Range("E21").Value = "  -5.60%  "
$ws.Range("E22").Value = "  -1.32%  "
Set-TextValue $ws "D23" "69.66"
$ws.Range("E23").Value = "  -2.58%  "
Set-TextValue $ws "D24" "265.52"
$ws.Range("E24").Value = "  -2.60%  "
Set-TextValue $ws "D25" "2.76"
$ws.Range("E25").Value = "  -1.30%  "
Set-TextValue $ws "D26" "0.180"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D27" "7.71"
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D28" "27.07"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  -0.10%  "
Set-TextValue $ws "D30" "0.108"
Set-TextValue $ws "D31" "6.49"
$ws.Range("E31").Value = "  +1.55%  "
Set-TextValue $ws "D32" "10.36"
$ws.Range("E32").Value = "  -3.83%  "
Set-TextValue $ws "D33" "36.57"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("E34").Value = "  +17.07%  "
Set-TextValue $ws "D35" "50.97"
$ws.Range("E35").Value = "  -4.07%  "
Set-TextValue $ws "D36" "0.0444"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -4.05%  "
$ws.Range("E39").Value = "  -2.28%  "
Set-TextValue $ws "D40" "17.99"
$ws.Range("E40").Value = "  -4.55%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  -0.75%  "
Set-TextValue $ws "D43" "23.06"
$ws.Range("E43").Value = "  -2.33%  "
Set-TextValue $ws "D44" "123.26"
$ws.Range("E44").Value = "  +9.04%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "2.131.53"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("E48").Value = "  -5.32%  "
$ws.Range("D49").Value = "3.312.20"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("E50").Value = "  +1.20%  "
Set-TextValue $ws "D51" "0.0332"
$ws.Range("E51").Value = "  -1.27%  "
